{"js": "// edit.js - Office.js (Word JavaScript API) script\n//\n// Applies two textual edits to the document:\n//\n// 1) The paragraph describing the convolution kernel used to contain two\n//    separate runs (split around a leftover \"_GoBack\" bookmark). The two\n//    runs are merged back into a single run/Text element with no\n//    intervening bookmark.\n//\n// 2) The SVM threshold value \"-1.5\" is changed to \"-2.5\" (i.e. the digit\n//    right after the minus sign changes from 1 to 2). Word records the\n//    edit point with a \"_GoBack\" bookmark placed right after the edited\n//    digit, which splits the run in two (\"-2\" and \".5\u662fSVM...\").\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Change 1: merge the \"\u7a0b\u5e8f\u4e2d\u4f7f\u7528\u7684\u5377\u79ef\u6838\u662f\" / \"\u7528PCAnet...\" runs into one\n// run and drop the bookmark that used to sit between them.\n// ---------------------------------------------------------------------\nconst mergedText =\n  \"\u7a0b\u5e8f\u4e2d\u4f7f\u7528\u7684\u5377\u79ef\u6838\u662f\u7528PCAnet\u5bf9FERET\u4eba\u8138\u6570\u636e\u96c6\u8bad\u7ec3\u5f97\u5230\uff0c\u53c2\u8003PCANet: A Simple Deep Learning Baseline for Image Classi\ufb01cation\u8fd9\u7bc7\u8bba\u6587\uff0c\u63d0\u53d6\u4eba\u8138\u7684\u7279\u5f81\uff0c\u5e76\u5c06\u4e24\u5f20\u8138\u4e4b\u95f4\u7684\u7279\u5f81\u5dee\u5f02\u7528\u5361\u65b9\u7edf\u8ba1\u4f5c\u4e3aSVM\u7684\u8f93\u5165\u6837\u672c\uff0c\u8fdb\u884c\u4e8c\u5206\u7c7b\u8bad\u7ec3\uff0c\u6807\u7b7e1\u8868\u793a\u4e24\u5f20\u56fe\u4ee3\u8868\u540c\u4e00\u4e2a\u4eba\uff0c\u6807\u7b7e0\u8868\u793a\u4e24\u5f20\u56fe\u4ee3\u8868\u4e0d\u540c\u4eba\uff0c\u751f\u6210\u4e00\u4e2a\u9a8c\u8bc1\u5668\u3002\u63d0\u53d6\u89c6\u9891\u56fe\u50cf\u4e2d\u68c0\u6d4b\u5230\u7684\u4eba\u8138\u7279\u5f81\u4e0e\u6a21\u677f\u7684\u4eba\u8138\u7279\u5f81\u6c42\u5361\u65b9\u7edf\u8ba1\u7684\u7279\u5f81\u5411\u91cf\uff0c\u9001\u5165\u8bad\u7ec3\u597d\u7684SVM\u505a\u5224\u522b\uff0c\u5f53\u5224\u5b9a\u7ed3\u679c\u4e3a1\u65f6\uff0c\u4e14\u7f6e\u4fe1\u5ea6\u9ad8\u4e8e\u6240\u8bbe\u9608\u503c\uff0c\u8ba4\u5b9a\u68c0\u6d4b\u5230\u4eba\u8138\u548c\u6a21\u677f\u56fe\u7247\u5c5e\u4e8e\u540c\u4e00\u4e2a\u4eba\u8138\u3002\";\n\nconst mergeResults = body.search(mergedText, { matchCase: true });\nmergeResults.load(\"text\");\nawait context.sync();\n\nconst mergedRange = mergeResults.items[0];\n// Re-inserting the identical text across the run/bookmark boundary collapses\n// everything (both runs + the bookmark) into a single fresh run.\nmergedRange.insertText(mergedText, Word.InsertLocation.replace);\nawait context.sync();\n\n// The stray bookmark ends up re-anchored at the end of the merged run;\n// the target document has no bookmark here at all, so drop it. (A new one\n// gets added below, at the second edit's location.)\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Change 2: \"-1.5\" -> \"-2.5\" (digit 1 -> 2), leaving the freshly-moved\n// \"_GoBack\" bookmark right after the new \"2\" (i.e. right before the \".\").\n// ---------------------------------------------------------------------\nconst thresholdResults = body.search(\"-1.5\", { matchCase: true });\nthresholdResults.load(\"text\");\nawait context.sync();\nconst thresholdRange = thresholdResults.items[0];\n\n// Scope a second search to just that match so we find the right \"1\".\nconst digitResults = thresholdRange.search(\"1\", { matchCase: true });\ndigitResults.load(\"text\");\nawait context.sync();\nconst digitRange = digitResults.items[0];\n\n// Zero-width point right after the digit (before the \".\") is where the new\n// \"_GoBack\" bookmark belongs, splitting the run in two.\nconst splitPoint = digitRange.getRange(Word.RangeLocation.after);\nsplitPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n\ndigitRange.insertText(\"2\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# edit.ps1 - Word COM interop script\n#\n# Applies two textual edits to the document:\n#\n# 1) The paragraph describing the convolution kernel used to contain two\n#    separate runs (split around a leftover \"_GoBack\" bookmark). The two\n#    runs are merged back into a single run/Text element with no\n#    intervening bookmark.\n#\n# 2) The SVM threshold value \"-1.5\" is changed to \"-2.5\" (i.e. the digit\n#    right after the minus sign changes from 1 to 2). Word records the\n#    edit point with a \"_GoBack\" bookmark placed right after the edited\n#    digit, which splits the run in two (\"-2\" and \".5\u662fSVM...\").\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Change 1: merge the \"\u7a0b\u5e8f\u4e2d\u4f7f\u7528\u7684\u5377\u79ef\u6838\u662f\" / \"\u7528PCAnet...\" runs into one\n# run and drop the bookmark that used to sit between them.\n# ---------------------------------------------------------------------\n$mergedText = \"\u7a0b\u5e8f\u4e2d\u4f7f\u7528\u7684\u5377\u79ef\u6838\u662f\u7528PCAnet\u5bf9FERET\u4eba\u8138\u6570\u636e\u96c6\u8bad\u7ec3\u5f97\u5230\uff0c\u53c2\u8003PCANet: A Simple Deep Learning Baseline for Image Classi\ufb01cation\u8fd9\u7bc7\u8bba\u6587\uff0c\u63d0\u53d6\u4eba\u8138\u7684\u7279\u5f81\uff0c\u5e76\u5c06\u4e24\u5f20\u8138\u4e4b\u95f4\u7684\u7279\u5f81\u5dee\u5f02\u7528\u5361\u65b9\u7edf\u8ba1\u4f5c\u4e3aSVM\u7684\u8f93\u5165\u6837\u672c\uff0c\u8fdb\u884c\u4e8c\u5206\u7c7b\u8bad\u7ec3\uff0c\u6807\u7b7e1\u8868\u793a\u4e24\u5f20\u56fe\u4ee3\u8868\u540c\u4e00\u4e2a\u4eba\uff0c\u6807\u7b7e0\u8868\u793a\u4e24\u5f20\u56fe\u4ee3\u8868\u4e0d\u540c\u4eba\uff0c\u751f\u6210\u4e00\u4e2a\u9a8c\u8bc1\u5668\u3002\u63d0\u53d6\u89c6\u9891\u56fe\u50cf\u4e2d\u68c0\u6d4b\u5230\u7684\u4eba\u8138\u7279\u5f81\u4e0e\u6a21\u677f\u7684\u4eba\u8138\u7279\u5f81\u6c42\u5361\u65b9\u7edf\u8ba1\u7684\u7279\u5f81\u5411\u91cf\uff0c\u9001\u5165\u8bad\u7ec3\u597d\u7684SVM\u505a\u5224\u522b\uff0c\u5f53\u5224\u5b9a\u7ed3\u679c\u4e3a1\u65f6\uff0c\u4e14\u7f6e\u4fe1\u5ea6\u9ad8\u4e8e\u6240\u8bbe\u9608\u503c\uff0c\u8ba4\u5b9a\u68c0\u6d4b\u5230\u4eba\u8138\u548c\u6a21\u677f\u56fe\u7247\u5c5e\u4e8e\u540c\u4e00\u4e2a\u4eba\u8138\u3002\"\n\n$find1 = $d.Content.Find\n$find1.Execute($mergedText, $false, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2)\n\n# ---------------------------------------------------------------------\n# Change 2: \"-1.5\" -> \"-2.5\" (digit 1 -> 2), leaving the freshly-moved\n# \"_GoBack\" bookmark right after the new \"2\" (i.e. right before the \".\").\n# ---------------------------------------------------------------------\n$r = $d.Content\n$find2 = $r.Find\n$find2.Text = \"-1.5\"\n$find2.Execute() | Out-Null\n\n$digitStart = $r.Start + 1     # the \"1\" in \"-1.5\"\n$splitPos = $r.Start + 2       # position between \"-1\" and \".5\"\n\n$rDigit = $d.Range($digitStart, $digitStart + 1)\n$rDigit.Text = \"2\"\n\n$rBookmark = $d.Range($splitPos, $splitPos)\n$d.Bookmarks.Add(\"_GoBack\", $rBookmark)\n"}
